# Update the CDA Logical model metadata for TEL (ST.r2b refresh):
#  - bump the IG Version string
#  - bump the publication Date
#  - insert a new "Jurisdiction" property/value row right after "Contact"
#    (pushing Description and everything below it down by one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Version (row 3, column B)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# 2) Date (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# 3) Insert a new row right below "Contact" (row 10) for "Jurisdiction"
$ws.Rows.Item(11).Insert()

# Match the formatting used by the other Property/Value rows (copy it down
# from the row that is now directly below the freshly inserted blank row).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
